$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colIndices = @(2, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13)

$data = @(
    @(13.50780914192253, 7.338121875690404, 13.33199612275783, 42.1540485701454, 50.32554851327314, 19.69306916641901, 32.18559989136059, 10.79730362532782, 11.19219560866518, 10.87546744287762, 15.45738209806132),
    @(13.42392698727772, 7.340436505544128, 13.34095018146238, 42.21097680419296, 50.34153512715321, 19.73600097076619, 32.28020930317126, 10.80322519421907, 11.00023668664309, 10.86658857259607, 15.44423674448167),
    @(13.37456945563479, 7.342211484088372, 13.34685289445694, 42.25421348363957, 50.36354498562978, 19.76544398661215, 32.3427900128795, 10.80709166478146, 10.88310061761984, 10.86275973212655, 15.43831145375423),
    @(13.35501303802626, 7.34302403596577, 13.34936029957271, 42.2739129238283, 50.37557708962757, 19.77821684413703, 32.36942121225089, 10.8087254093183, 10.83561450347782, 10.86160934632285, 15.43643936744067),
    @(13.35179985538147, 7.34316435701419, 13.34978281856, 42.27730955054174, 50.3777598719762, 19.78038453248482, 32.37391149382037, 10.80900020615736, 10.8277463379041, 10.8614431319864, 15.43616135523533),
    @(13.37430343253134, 7.342222080770841, 13.34688629693153, 42.25447073860976, 50.36369485897184, 19.76561311039325, 32.34314459912225, 10.80711346249123, 10.88245911406602, 10.8627425556979, 15.43828400588115),
    @(13.47845177988268, 7.33884669388628, 13.33499956222902, 42.17195700495407, 50.32852828681564, 19.70723190030655, 32.21728926916649, 10.79929760428302, 11.12589594299744, 10.87207034212793, 15.45240597561257),
    @(13.69891585913828, 7.335023170556464, 13.31489388731634, 42.07596115163273, 50.35640364301604, 19.61723129133325, 32.00611585293647, 10.78579419989198, 11.606231535824, 10.90315130438847, 15.49699147643647),
    @(13.86966529939782, 7.333903353939604, 13.30206422724924, 42.04565600431941, 50.4359313450564, 19.56606632663891, 31.87268808333333, 10.7769762500267, 11.95725140763927, 10.93365610604467, 15.5398527283433),
    @(13.94899900686952, 7.333757573533589, 13.29664700039229, 42.04061277815598, 50.48489786389491, 19.54604376333945, 31.81670545286131, 10.77320242585289, 12.11578565075243, 10.94916783768992, 15.56149921620858),
    @(13.9792574998543, 7.333754382155886, 13.29465571652451, 42.03995976892127, 50.50527305285672, 19.53892980009214, 31.79618435934538, 10.77180739214864, 12.17559463638238, 10.95527384207127, 15.57000079703749),
    @(13.97273151474921, 7.333752760641597, 13.2950819045303, 42.04004452615199, 50.50080350683331, 19.54044109252607, 31.80057376942067, 10.77210632596932, 12.16272460553644, 10.95394853509138, 15.56815635911565),
    @(13.95148417859224, 7.333756270161943, 13.29648197261436, 42.04053387250989, 50.48653745776025, 19.54544910880804, 31.81500357343385, 10.77308697423572, 12.12071100287791, 10.94966554994633, 15.56219256960097),
    @(13.93849711343632, 7.33376518515476, 13.29734737725125, 42.04099725250686, 50.47803753080906, 19.54857764109424, 31.82393059340971, 10.77369207795954, 12.09494549994194, 10.94707222339933, 15.55857909824607),
    @(13.86451216206462, 7.333920174556679, 13.30242667467814, 42.04616149488412, 50.43298811038183, 19.56744033745207, 31.87644157825079, 10.7772276471159, 11.94686259259386, 10.93267502053167, 15.53848094146588),
    @(13.81953411902032, 7.334108193325973, 13.30564987559588, 42.0515687157755, 50.40862335735908, 19.57984539394979, 31.9098631020167, 10.77945734824166, 11.85568160510682, 10.92425944065688, 15.52669861330349),
    @(13.793821424135, 7.334250580037073, 13.30754322768951, 42.05550178346535, 50.39581365815396, 19.58728658742781, 31.92952998941528, 10.78076217566026, 11.80313220500666, 10.91957310579989, 15.52012425667114),
    @(13.78514326153787, 7.334304679444551, 13.30819106492721, 42.05697480215071, 50.39168351144936, 19.58985861222097, 31.93626505531125, 10.7812078122733, 11.78532379318037, 10.91801294728338, 15.51793320836292),
    @(13.82430596781225, 7.334084636368823, 13.30530267835753, 42.05090793382531, 50.41109243058969, 19.57849316755049, 31.90625939844259, 10.77921767922298, 11.8653992336824, 10.92513936787393, 15.52793193295906),
    @(13.95771933624458, 7.333753829944192, 13.29606910875269, 42.04035603869072, 50.49067806559429, 19.54396542468746, 31.81074677777295, 10.77279801142047, 12.13305797056808, 10.95091729409791, 15.56393605109809),
    @(14.04616445847061, 7.333840667593829, 13.29038468926479, 42.04078438082464, 50.55336954160259, 19.52412840636001, 31.75227767792825, 10.76880069445235, 12.3066542518692, 10.96911565958684, 15.58923959689114),
    @(13.99885243406554, 7.33376668463165, 13.29338657379081, 42.0398858878747, 50.51893552194043, 19.53446598705193, 31.78312182023827, 10.77091603232883, 12.21414370343761, 10.95928029237716, 15.57557394290199),
    @(13.8221481592094, 7.334095179655943, 13.30545952073745, 42.05120410535352, 50.40997243058813, 19.57910354545619, 31.90788722358999, 10.779325962119, 11.86100628864924, 10.92474107940728, 15.52737372750199),
    @(13.63765016994557, 7.335759783825177, 13.31999116190766, 42.0948718645155, 50.33848699148133, 19.63895401365242, 32.05942948380078, 10.78925291199823, 11.47635840382413, 10.89338621824158, 15.483140474173)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $colIndices.Length; $c++) {
        $ws.Cells.Item($r + 2, $colIndices[$c]).Value = $rowValues[$c]
    }
}
